# ProstateBrachyQA Log.xlsx edit
# - Drop the "Test" column/series (column C) entirely: it was deleted from
#   the sheet, so the chart series that plotted it is removed too.
# - Update the Weight value for the 3rd data point (B4): 1200 -> 1300.
# - Selection moves to A4:B5 (active cell A5).
# - Chart gets an (auto) title turned on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phantom")

# --- Worksheet data edits -------------------------------------------------

# Weight column: 3rd point changes from 1200 to 1300.
$ws.Range("B4").Value = 1300

# "Test" column (C) is removed: clear its header (C1) and only data value (C4).
$ws.Range("C1").ClearContents() | Out-Null
$ws.Range("C4").ClearContents() | Out-Null

# Selection now covers A4:B5.
$ws.Range("A4:B5").Select() | Out-Null

# --- Chart edits -----------------------------------------------------------

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Remove the second series ("Test"), which plotted the now-empty column C.
$chart.SeriesCollection().Item(2).Delete() | Out-Null

# Turn on the (automatic) chart title.
$chart.HasTitle = $true
